$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("costumes")

# Row 1 - headers: add new "accessories" column at H, shift old H (feet) to I
$ws.Range("H1").Value = "accessories"
$ws.Range("I1").Value = "feet"
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Columns.Item(9).ColumnWidth = 9.7265625

# Row 2 - Daria Morgendorffer
$ws.Range("D2").Value = "brown hair"
$ws.Range("E2").Value = "green blazer, green jacket, orange shirt, orange top, orange blouse"
$ws.Range("F2").Value = "black skirt"
$ws.Range("H2").Value = "round glasses"
$ws.Range("I2").Value = "black boots, combat boots, black combat boots"

# Row 3 - Jane Lane
$ws.Range("D3").Value = "black hair, black bob"
$ws.Range("E3").Value = "red blazer, red jacket, black shirt, black blouse, black v-neck"
$ws.Range("F3").Value = "grey shorts, black shorts"
$ws.Range("H3").Value = "earrings"
$ws.Range("I3").Value = "black tights, black boots, black knee-high boots"

# Row 4 - Quinn Morgendorffer
$ws.Range("A4").Value = "Daria"
$ws.Range("B4").Value = "Quinn Morgendorffer"
$ws.Range("C4").Value = "https://th.bing.com/th/id/R.e65c73bba6f6ae5d6f882f7932ce2d4c?rik=7UN34tNBhRF77g&riu=http%3a%2f%2fimg2.wikia.nocookie.net%2f__cb20130604070322%2fdaria%2fimages%2f8%2f8b%2fImg-thing.jpg&ehk=Df7ytuLZKvLFc%2bfDnDQdfQCoNrc%2f%2btcAowq7t5B2UtE%3d&risl=&pid=ImgRaw&r=0"
$ws.Range("D4").Value = "red hair, bangs, redhead"
$ws.Range("E4").Value = "pink shirt, smiley face tee"
$ws.Range("F4").Value = "flare jeans, blue jeans, flared pants, blue pants"
$ws.Range("H4").Value = "brown belt"
$ws.Range("I4").Value = "black boots, black shoes"

# Row 5 - Trent Lane
$ws.Range("A5").Value = "Daria"
$ws.Range("B5").Value = "Trent Lane"
$ws.Range("C5").Value = "https://static.wikia.nocookie.net/daria/images/f/f3/Trentcropped.png/revision/latest?cb=20200815062032"
$ws.Range("D5").Value = "spiky black hair, goatee"
$ws.Range("E5").Value = "green shirt, green tee"
$ws.Range("H5").Value = "necklace, cuff bracelet, rings, earrings"

# Row 6 - Brittany Taylor
$ws.Range("A6").Value = "Daria"
$ws.Range("B6").Value = "Brittany Taylor"
$ws.Range("C6").Value = "https://static.miraheze.org/dariawikiwiki/c/c5/Brittany.gif"
$ws.Range("D6").Value = "blonde pigtails, blond hair"
$ws.Range("E6").Value = "blue shirt"
$ws.Range("F6").Value = "yellow skirt, blue skirt, pleated skirt"
$ws.Range("H6").Value = "pom-poms"
$ws.Range("I6").Value = "yellow socks, blue sneakers, blue shoes"

# Row 7 - Kevin Thompson
$ws.Range("A7").Value = "Daria"
$ws.Range("B7").Value = "Kevin Thompson"
$ws.Range("C7").Value = "https://th.bing.com/th/id/R.58bd988ef1a6fdb387eb2c760efc707f?rik=qbnGZ8NksR7AjA&riu=http%3a%2f%2fvignette4.wikia.nocookie.net%2fdaria%2fimages%2f1%2f14%2fKevin_Thompson.gif%2frevision%2flatest%3fcb%3d20140902121519&ehk=zzgrLHBjr0sgpwmYGZ%2fVDP43hK2tTkgKGckfZC7drhc%3d&risl=&pid=ImgRaw&r=0"
$ws.Range("D7").Value = "black hair"
$ws.Range("E7").Value = "yellow football jersey, yellow shirt"
$ws.Range("F7").Value = "blue pants, jeans, blue leggings"
$ws.Range("H7").Value = "neck cushion, football equipment, gloves"
$ws.Range("I7").Value = "cleats, white sneakers"

# Row 8 - Jodie Landon
$ws.Range("A8").Value = "Daria"
$ws.Range("B8").Value = "Jodie Landon"
$ws.Range("C8").Value = "https://static.tvtropes.org/pmwiki/pub/images/jodie_8099.gif"
$ws.Range("D8").Value = "black hair, black braids"
$ws.Range("E8").Value = "pink shirt, pink blouse"
$ws.Range("F8").Value = "grey skirt, white skirt, miniskirt"
$ws.Range("H8").Value = "books"
$ws.Range("I8").Value = "black shoes, black mary janes, black loafers"

$ws.Range("A9").Select() | Out-Null
